$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8099.407
$ws.Range("I76").Value = 8588.632
$ws.Range("J76").Value = 6937.5
$ws.Range("K76").Value = 8588.632
$ws.Range("L76").Value = 6937.5
$ws.Range("M76").Value = -8273.632
$ws.Range("N76").Value = -7567.5

$ws.Range("H79").Value = 8099.407
$ws.Range("I79").Value = 8588.632
$ws.Range("J79").Value = 6937.5
$ws.Range("K79").Value = 8588.632
$ws.Range("L79").Value = 6937.5
$ws.Range("M79").Value = -7496.632
$ws.Range("N79").Value = -9121.5

$ws.Range("H98").Value = 1146.1666
$ws.Range("I98").Value = 1168.4482
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 1168.4482
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 329.5518

$ws.Range("H100").Value = 3915.2285
$ws.Range("I100").Value = 1892.4286
$ws.Range("J100").Value = 6949.4287
$ws.Range("K100").Value = 1892.4286
$ws.Range("L100").Value = 6949.4287
$ws.Range("M100").Value = -1351.4286
$ws.Range("N100").Value = -8031.4287

$ws.Range("H107").Value = 798.6316
$ws.Range("I107").Value = 833.7646999999999
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 833.7646999999999
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1086.2353

$ws.Range("H111").Value = 522.8570999999999
$ws.Range("I111").Value = 200
$ws.Range("J111").Value = 652
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 1956
$ws.Range("M111").Value = 2467
$ws.Range("N111").Value = -8090

$ws.Range("H122").Value = 1146.1666
$ws.Range("I122").Value = 1168.4482
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 3505.3446
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -1055.3446

$ws.Range("H133").Value = 57802
$ws.Range("I133").Value = 40010
$ws.Range("J133").Value = 62250
$ws.Range("K133").Value = 40010
$ws.Range("L133").Value = 62250
$ws.Range("M133").Value = -34950
$ws.Range("N133").Value = -72370

$ws.Range("H136").Value = 38756
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 38756
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 38756
$ws.Range("N136").Value = -48956

$ws.Range("H137").Value = 441839.06
$ws.Range("I137").Value = 3689.476
$ws.Range("J137").Value = 1464188.1
$ws.Range("K137").Value = 11068.428
$ws.Range("L137").Value = 4392564.300000001
$ws.Range("M137").Value = -8518.428
$ws.Range("N137").Value = -4397664.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23546.314
$ws.Range("I32").Value = 14092.392
$ws.Range("J32").Value = 41666.332
$ws.Range("K32").Value = 14092.392
$ws.Range("L32").Value = 41666.332
$ws.Range("M32").Value = -13805.392
$ws.Range("N32").Value = -42240.332

$ws.Range("H132").Value = 2325.7195
$ws.Range("I132").Value = 2122.9016
$ws.Range("J132").Value = 2914.8572
$ws.Range("K132").Value = 6368.7048
$ws.Range("L132").Value = 8744.571599999999
$ws.Range("M132").Value = -3838.7048
$ws.Range("N132").Value = -13804.5716

$ws.Range("H134").Value = 34479.8
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 34479.8
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 34479.8
$ws.Range("N134").Value = -44619.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 275.53125
$ws.Range("I80").Value = 106.3
$ws.Range("J80").Value = 352.45456
$ws.Range("K80").Value = 106.3
$ws.Range("L80").Value = 352.45456
$ws.Range("M80").Value = 891.7
$ws.Range("N80").Value = -2348.45456

$ws.Range("H83").Value = 275.53125
$ws.Range("I83").Value = 106.3
$ws.Range("J83").Value = 352.45456
$ws.Range("K83").Value = 531.5
$ws.Range("L83").Value = 1762.2728
$ws.Range("M83").Value = 4460.5
$ws.Range("N83").Value = -11746.2728

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 415.77777
$ws.Range("I107").Value = 288.94116
$ws.Range("J107").Value = 631.4
$ws.Range("K107").Value = 288.94116
$ws.Range("L107").Value = 631.4
$ws.Range("M107").Value = 1631.05884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 690.59576
$ws.Range("I5").Value = 476.96875
$ws.Range("J5").Value = 1146.3334
$ws.Range("K5").Value = 1430.90625
$ws.Range("L5").Value = 3439.0002
$ws.Range("M5").Value = -1318.90625
$ws.Range("N5").Value = -3663.0002

$ws.Range("H62").Value = 3400
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3633.3333
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 10899.9999
$ws.Range("M62").Value = -5314
$ws.Range("N62").Value = -12271.9999

$ws.Range("H65").Value = 3400
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3633.3333
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 32699.9997
$ws.Range("M65").Value = -14568
$ws.Range("N65").Value = -39563.9997

$ws.Range("H68").Value = 837.6
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 922
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2766
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -4388

$ws.Range("H70").Value = 4988.8887
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 5985.7144
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 17957.1432
$ws.Range("M70").Value = -4185
$ws.Range("N70").Value = -18587.1432

$ws.Range("H71").Value = 837.6
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 922
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 8298
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -16410

$ws.Range("H73").Value = 4988.8887
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 5985.7144
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 17957.1432
$ws.Range("M73").Value = -3408
$ws.Range("N73").Value = -20141.1432

$ws.Range("H97").Value = 346.1
$ws.Range("I97").Value = 228.6
$ws.Range("J97").Value = 463.6
$ws.Range("K97").Value = 685.8
$ws.Range("L97").Value = 1390.8
$ws.Range("M97").Value = -189.8
$ws.Range("N97").Value = -2382.8

$ws.Range("H98").Value = 4388.6665
$ws.Range("I98").Value = 234.5
$ws.Range("J98").Value = 5219.5
$ws.Range("K98").Value = 703.5
$ws.Range("L98").Value = 15658.5
$ws.Range("M98").Value = 794.5
$ws.Range("N98").Value = -18654.5

$ws.Range("H122").Value = 559.73334
$ws.Range("I122").Value = 422.84616
$ws.Range("J122").Value = 1449.5
$ws.Range("K122").Value = 3805.61544
$ws.Range("L122").Value = 13045.5
$ws.Range("M122").Value = -1355.61544

$ws.Range("H135").Value = 690.59576
$ws.Range("I135").Value = 476.96875
$ws.Range("J135").Value = 1146.3334
$ws.Range("K135").Value = 4292.71875
$ws.Range("L135").Value = 10317.0006
$ws.Range("M135").Value = -1757.71875
$ws.Range("N135").Value = -15387.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.3125
$ws.Range("I2").Value = 15.5
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 15.5
$ws.Range("L2").Value = 55
$ws.Range("M2").Value = 97.5
$ws.Range("N2").Value = -281

$ws.Range("H97").Value = 1208.2307
$ws.Range("I97").Value = 994.0625
$ws.Range("J97").Value = 1550.9
$ws.Range("K97").Value = 994.0625
$ws.Range("L97").Value = 1550.9
$ws.Range("M97").Value = -498.0625
$ws.Range("N97").Value = -2542.9

$ws.Range("H113").Value = 5221.0454
$ws.Range("I113").Value = 7688.5
$ws.Range("J113").Value = 903
$ws.Range("K113").Value = 7688.5
$ws.Range("L113").Value = 903
$ws.Range("M113").Value = -5518.5
$ws.Range("N113").Value = -5243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2027.3572
$ws.Range("I61").Value = 2053.5
$ws.Range("J61").Value = 1870.5
$ws.Range("K61").Value = 2053.5
$ws.Range("L61").Value = 1870.5
$ws.Range("M61").Value = -1851.5
$ws.Range("N61").Value = -2274.5

$ws.Range("H100").Value = 31253394
$ws.Range("I100").Value = 4481.6113
$ws.Range("J100").Value = 71430570
$ws.Range("K100").Value = 4481.6113
$ws.Range("L100").Value = 71430570
$ws.Range("M100").Value = -3940.6113
$ws.Range("N100").Value = -71431652

$ws.Range("H113").Value = 2027.3572
$ws.Range("I113").Value = 2053.5
$ws.Range("J113").Value = 1870.5
$ws.Range("K113").Value = 2053.5
$ws.Range("L113").Value = 1870.5
$ws.Range("M113").Value = 116.5
$ws.Range("N113").Value = -6210.5

$ws.Range("H114").Value = 24000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 24000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 24000
$ws.Range("N114").Value = -32678

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H123").Value = 13429
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 13429
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 13429
$ws.Range("N123").Value = -23229

$ws.Range("H132").Value = 2426.1162
$ws.Range("I132").Value = 1583.0358
$ws.Range("J132").Value = 3999.8667
$ws.Range("K132").Value = 4749.107400000001
$ws.Range("L132").Value = 11999.6001
$ws.Range("M132").Value = -2219.107400000001
$ws.Range("N132").Value = -17059.6001

$ws.Range("H135").Value = 56270
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 56270
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 56270
$ws.Range("N135").Value = -66410
